$wb = $excel.ActiveWorkbook

# --- Swiss sheet: remove the "XLM800-STI" / "XLM800-Zetfas" rows (old rows 9-10) ---
$wsSwiss = $wb.Worksheets.Item("Swiss")
$wsSwiss.Activate()
$wsSwiss.Rows("9:10").Select()
$wsSwiss.Rows("9:10").Delete()

# --- Portugal sheet: same row removal, then leave it as the active sheet/cell ---
$wsPortugal = $wb.Worksheets.Item("Portugal")
$wsPortugal.Activate()
$wsPortugal.Rows("9:10").Select()
$wsPortugal.Rows("9:10").Delete()
$wsPortugal.Range("A8").Select()
